$d = $word.ActiveDocument

# Title on page 1: "Tic-Tac-Toe" -> "Connect Four"
$d.Content.Find.Execute("Tic-Tac-Toe", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Connect Four", 2)

# Overview paragraph: "we've chosen to create" -> "we will create"
# (curly apostrophe, U+2019, matches the source run)
$apos = [char]0x2019
$d.Content.Find.Execute("we${apos}ve chosen to create", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "we will create", 2)
